# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") for rows 2-36, replacing the old Strike# values
$newK = @{
    2  = 6
    3  = 4
    4  = 5
    5  = 5
    6  = 8
    7  = 5
    8  = 2
    9  = 7
    10 = 7
    11 = 6
    12 = 8
    13 = 6
    14 = 4
    15 = 5
    16 = 8
    17 = 3
    18 = 7
    19 = 5
    20 = 8
    21 = 6
    22 = 10
    23 = 4
    24 = 8
    25 = 3
    26 = 5
    27 = 3
    28 = 4
    29 = 4
    30 = 4
    31 = 6
    32 = 4
    33 = 0
    34 = 1
    35 = 3
    36 = 4
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
